$d = $word.ActiveDocument

# Namespace + shared formatting fragments used to build the new
# ("Electricity Department") paragraph's raw OOXML.
$wns  = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$rPr  = "<w:rPr><w:rFonts w:ascii=`"Bookman Old Style`" w:hAnsi=`"Bookman Old Style`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"18`"/><w:szCs w:val=`"20`"/></w:rPr>"
$pPr  = "<w:pPr><w:ind w:left=`"113`" w:right=`"113`"/><w:jc w:val=`"center`"/>$rPr</w:pPr>"
$newParaXml = "<w:p $wns>$pPr<w:r>$rPr<w:t>Electricity Department</w:t></w:r></w:p>"

$searchText = "Circle-I Electricity Department"

$rng = $d.Content
$rng.Find.Forward = $true
$rng.Find.Wrap = 0

while ($rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    # Shrink the run's visible text down to "Circle-I " (trailing space
    # kept, hence xml:space="preserve" is emitted automatically) while
    # leaving the paragraph's own <w:p>/<w:pPr> markup untouched.
    $rng.Text = "Circle-I "

    # Insert a brand-new paragraph right after it, carrying
    # "Electricity Department" with matching paragraph/run formatting.
    $insertPoint = $d.Range($rng.End, $rng.End)
    $insertPoint.InsertXML($newParaXml)

    # Continue searching after the content we just touched.
    $rng.Start = $insertPoint.End
    $rng.End = $d.Content.End
}
